$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-detected as numbers by Excel (these are text-formatted price cells).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'

# Apply the updated values.
$ws.Range('D2').Value = '70.663.84'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '3.629.20'
$ws.Range('E3').Value = '  +3.67%  '
$ws.Range('D5').Value = '607.51'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = '198.93'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.220'
$ws.Range('E9').Value = '  +9.85%  '
$ws.Range('D10').Value = '0.646'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('D11').Value = '53.78'
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '4.204.76'
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('D15').Value = '684.24'
$ws.Range('E15').Value = '  +15.01%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.698.60'
$ws.Range('E16').Value = '  +5.75%  '
$ws.Range('D17').Value = '12.94'
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.818.54'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').Value = '18.73'
$ws.Range('E22').Value = '  +3.02%  '
$ws.Range('D23').Value = '5.37'
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('D24').Value = '105.23'
$ws.Range('E24').Value = '  +3.57%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  -4.33%  '
$ws.Range('D27').Value = '10.48'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = '9.96'
$ws.Range('E28').Value = '  +4.97%  '
$ws.Range('D29').Value = '34.20'
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('D30').Value = '4.56'
$ws.Range('E30').Value = '  +6.25%  '
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('D32').Value = '12.18'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = '63.27'
$ws.Range('D35').Value = '3.960.53'
$ws.Range('E35').Value = '  +6.04%  '
$ws.Range('D36').Value = '0.0₃0867'
$ws.Range('E36').Value = '  +6.73%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').Value = '36.64'
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('D40').Value = '0.387'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '3.55'
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '498.14'
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('D43').Value = '0.136'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('D44').Value = '3.07'
$ws.Range('E44').Value = '  +9.34%  '
$ws.Range('D45').Value = '0.0457'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').Value = '3.50'
$ws.Range('E46').Value = '  +6.31%  '
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('E51').Value = '  +1.85%  '

# Restore default (unstyled) formatting now that the text values are locked in,
# so these cells keep matching the original un-styled inline-string cells.
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
